# "Worked on reroute request functionality"
#
# The "Reroute Request" sheet (sheet1) grows from a 9-column x 2-row
# reference table into a 15-column x 4-row table: three new columns are
# inserted at the front (Location Type / OrderId / Way Bill #), two more
# are inserted between "Country" and "Residential Location"
# (Earliest Drop-Off / Latest Drop-Off), and a trailing "orderReferenceID"
# column is appended. Two more sample rows are added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reroute Request")

# --- Make room: shift the existing columns to their new homes -------------
# old A:F (Location Name..Country) -> new D:I
# old G:I (Residential/Commercial/Warehouse Location) -> new L:N (after
# the two new Drop-Off columns land on J:K)
$ws.Columns("A:C").Insert()
$ws.Columns("J:K").Insert()

# --- Row 1: headers ---------------------------------------------------
$ws.Range("A1").Value = "Location Type"
$ws.Range("B1").Value = "OrderId"
$ws.Range("C1").Value = "Way Bill #"
$ws.Range("J1").Value = "Earliest Drop-Off"
$ws.Range("K1").Value = "Latest Drop-Off"
$ws.Range("O1").Value = "orderReferenceID"

# --- Row 2: existing sample row gains new columns ----------------------
$ws.Range("A2").Value = "Residential"
$ws.Range("B2").Formula = "'51487615"
$ws.Range("C2").Value = "CEV1002186"
$ws.Range("F2").Value = "LOS ANGELES"
$ws.Range("G2").Value = "CA"
$ws.Range("H2").Formula = "'90001"
$ws.Range("L2").Value = "Residential"
$ws.Range("M2").Value = "Commercial"
$ws.Range("N2").Value = "Amazon FBA Warehouse"

# --- Row 3: new sample row ----------------------------------------------
$ws.Range("A3").Value = "Commercial"
$ws.Range("B3").Formula = "'51488422"
$ws.Range("C3").Formula = "'66488699"
$ws.Range("D3").Value = "Test Location 2"
$ws.Range("E3").Value = "#123, Moody"
$ws.Range("F3").Value = "MOODY"
$ws.Range("G3").Value = "AL"
$ws.Range("H3").Formula = "'35004"
$ws.Range("I3").Value = "US"
$ws.Range("J3").NumberFormat = "h:mm AM/PM"
$ws.Range("J3").Value = 0.5
$ws.Range("K3").NumberFormat = "h:mm AM/PM"
$ws.Range("K3").Value = 2/3

# --- Row 4: new sample row ----------------------------------------------
$ws.Range("A4").Value = "Amazon FBA Warehouse"
$ws.Range("B4").Formula = "'51488421"
$ws.Range("C4").Formula = "'66488700"
$ws.Range("D4").Value = "Test Location 3"
$ws.Range("E4").Value = "#321, Los Angeles"
$ws.Range("F4").Value = "LOS ANGELES"
$ws.Range("G4").Value = "CA"
$ws.Range("H4").Formula = "'90001"
$ws.Range("I4").Value = "US"

# --- Column widths for the newly-inserted columns -----------------------
$ws.Columns("A").ColumnWidth = 133/6
$ws.Columns("B:C").ColumnWidth = 12.5
$ws.Columns("J:K").ColumnWidth = 91/6

# --- Page setup + view state ---------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.Range("D8").Select()
